# "viz y tablas update"
#
# Updates the "Ficha tecnica" (fact sheet) metadata table on the second
# worksheet:
#   - removes the "DIMENSION" / "Calidad" key-value row
#   - adds a new "TIPOIND" / "Resultados" key-value row
#   - adds a new "CITA" / "UMAD con base en SINADI - MSP" key-value row
#
# The metadata table lives in columns A (key) and B (value) of the
# "Ficha tecnica" worksheet, one key/value pair per row, starting at row 1.
# Before the edit it looks like:
#   1: (blank) | " "
#   2: DERECHO | Salud
#   3: DIMENSION | Calidad
#   4: CONINDICADOR | Razon de consultas no urgentes/ consultas urgentes
#   5: NOMINDICADOR | Razon Consultas No Urgentes/ Consultas Urgentes
#   6: DEFINICION | El indicador mide ...
#   7: CALCULO | Para cada anio calcular ...
#
# Removing row 3 (DIMENSION/Calidad) shifts rows 4-7 up by one, and the two
# new rows (TIPOIND/Resultados, CITA/UMAD con base en SINADI - MSP) are
# appended after what is left, giving 8 rows in total:
#   1: (blank) | " "
#   2: DERECHO | Salud
#   3: CONINDICADOR | Razon de consultas no urgentes/ consultas urgentes
#   4: NOMINDICADOR | Razon Consultas No Urgentes/ Consultas Urgentes
#   5: DEFINICION | El indicador mide ...
#   6: CALCULO | Para cada anio calcular ...
#   7: TIPOIND | Resultados
#   8: CITA | UMAD con base en SINADI - MSP

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Capture the existing key/value pairs for the rows that will shift up
# (rows 4-7) before overwriting anything.
$key4 = $ws.Cells.Item(4, 1).Value2
$val4 = $ws.Cells.Item(4, 2).Value2
$key5 = $ws.Cells.Item(5, 1).Value2
$val5 = $ws.Cells.Item(5, 2).Value2
$key6 = $ws.Cells.Item(6, 1).Value2
$val6 = $ws.Cells.Item(6, 2).Value2
$key7 = $ws.Cells.Item(7, 1).Value2
$val7 = $ws.Cells.Item(7, 2).Value2

# Shift rows 4-7 up into rows 3-6, effectively removing the
# "DIMENSION" / "Calidad" row that used to sit at row 3.
$ws.Cells.Item(3, 1).Value = $key4
$ws.Cells.Item(3, 2).Value = $val4
$ws.Cells.Item(4, 1).Value = $key5
$ws.Cells.Item(4, 2).Value = $val5
$ws.Cells.Item(5, 1).Value = $key6
$ws.Cells.Item(5, 2).Value = $val6
$ws.Cells.Item(6, 1).Value = $key7
$ws.Cells.Item(6, 2).Value = $val7

# Append the new metadata rows.
$ws.Cells.Item(7, 1).Value = "TIPOIND"
$ws.Cells.Item(7, 2).Value = "Resultados"
$ws.Cells.Item(8, 1).Value = "CITA"
$ws.Cells.Item(8, 2).Value = "UMAD con base en SINADI - MSP"
